# Auto-generated edit script to update cryptos list values
# matching the commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.012.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.815.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.62%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.614"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("E7").Value = "  +0.64%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.26"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -10.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.319"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.80%  "
$ws.Range("E10").Value = "  +0.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0999"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.074.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.811.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("E15").Value = "  +4.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.659"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.971.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "238.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.68%  "
$ws.Range("E21").Value = "  +1.01%  "
$ws.Range("E22").Value = "  +3.81%  "
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("E24").Value = "  +5.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.44%  "
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +32.89%  "
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.339.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E32").Value = "  +6.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("E35").Value = "  -3.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "93.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.681"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.35%  "
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("E40").Value = "  +5.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.310.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("E42").Value = "  +3.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.20%  "
$ws.Range("E44").Value = "  -3.95%  "
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0512"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.990.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.62%  "
$ws.Range("E51").Value = "  +5.35%  "
